$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F12").Value = '{"de": "Ausgewählte Atérien Fundstellen aus der ROAD Datenbank", "en": "Selected Aterian sites from the ROAD Database"}'
$ws.Range("F13").Value = '{"de": "Ausgewählte Micoquien Fundstellen aus der ROAD Datenbank", "en": "Selected Micoquian sites from the ROAD Database"}'
$ws.Range("F14").Value = '{"de": "Ausgewählte Gravettien Fundstellen aus der ROAD Datenbank", "en": "Selected Gravettian sites from the ROAD Database"}'
$ws.Range("F19").Value = '{"de": "Ausgewählte Still Bay Fundstellen aus der ROAD Datenbank", "en": "Selected Still Bay sites from the ROAD Database"}'
$ws.Range("F20").Value = '{"de": "Ausgewählte Howiesonspoort Fundstellen aus der ROAD Datenbank", "en": "Selected Howiesonspoort sites from the ROAD Database"}'
$ws.Range("F23").Value = '{"de": "Ausgewählte Châtelperronien Fundstellen aus der ROAD Datenbank", "en": "Selected Châtelperronian sites from the ROAD Database"}'
$ws.Range("F24").Value = '{"de": "Ausgewählte Ahmarien Fundstellen aus der ROAD Datenbank", "en": "Selected Ahmarian sites from the ROAD Database"}'
$ws.Range("F27").Value = '{"de": "Ausgewählte La Quina Fundstellen aus der ROAD Datenbank", "en": "Selected La Quina sites from the ROAD Database"}'
$ws.Range("F29").Value = '{"de": "Ausgewählte Clactionien Fundstellen aus der ROAD Datenbank", "en": "Selected Clactionian sites from the ROAD Database"}'
$ws.Range("F31").Value = '{"de": "Ausgewählte Fauresmith Fundstellen aus der ROAD Datenbank", "en": "Selected Fauresmith sites from the ROAD Database"}'
$ws.Range("F32").Value = '{"de": "Ausgewählte Emiran Fundstellen aus der ROAD Datenbank", "en": "Selected Emiran sites from the ROAD Database"}'
$ws.Range("F37").Value = '{"de": "Ausgewählte Szeletien Fundstellen aus der ROAD Datenbank", "en": "Selected Szeletian sites from the ROAD Database"}'

$ws.Range("F38").Select()
